$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of translation data to append below the existing table (rows 11-19)
$data = @(
    @("Graphics: ", "Graphics: ", "图像"),
    @("Volume:", "Volume:", "音量"),
    @("Language:", "Language:", "语言"),
    @("Settings", "Settings", "设置"),
    @("Deep Dive Descent", "Deep Dive Descent", "深浅迷航"),
    @("Language", "Language", "语言"),
    @("Return", "Return", "返回"),
    @("English", "English", "English"),
    @("Chinese", "Chinese", "中文")
)

$row = 11
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

# Update selection to match the new last cell used in the sheet
$ws.Range("C19").Select()
